$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) get re-styled from the custom
#    "Table_0" style to the built-in table style
#    {C3234168-431A-4032-BC5D-251C5E70CF5D}.
# ---------------------------------------------------------------------------
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle("{C3234168-431A-4032-BC5D-251C5E70CF5D}")
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme colour scheme switches from the "Integral" design's
#    "Red Violet" palette to the default Office palette.
# ---------------------------------------------------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
$officeRgb = @(
    0,         # dk1     000000
    16777215,  # lt1     FFFFFF
    6968388,   # dk2     44546A
    15132391,  # lt2     E7E6E6
    13998939,  # accent1 5B9BD5
    3243501,   # accent2 ED7D31
    10855845,  # accent3 A5A5A5
    49407,     # accent4 FFC000
    12874308,  # accent5 4472C4
    4697456,   # accent6 70AD47
    12673797,  # hlink   0563C1
    7491477    # folHlink 954F72
)
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeRgb[$i - 1]
}
